$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-10 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-11 Thursday", 2)

$d.Content.Find.Execute("482÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "537÷3=", 2)
$d.Content.Find.Execute("738÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "721÷7=", 2)
$d.Content.Find.Execute("741÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "827÷5=", 2)
$d.Content.Find.Execute("931÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "789÷2=", 2)
$d.Content.Find.Execute("893÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "192÷3=", 2)
$d.Content.Find.Execute("877÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "164÷7=", 2)
$d.Content.Find.Execute("501÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "986÷8=", 2)
$d.Content.Find.Execute("317÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "236÷9=", 2)
$d.Content.Find.Execute("351÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "246÷6=", 2)
$d.Content.Find.Execute("419÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "839÷7=", 2)
$d.Content.Find.Execute("882÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "517÷4=", 2)
$d.Content.Find.Execute("846÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "109÷5=", 2)
$d.Content.Find.Execute("338÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "988÷9=", 2)
$d.Content.Find.Execute("895÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "627÷6=", 2)
$d.Content.Find.Execute("781÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "896÷9=", 2)
$d.Content.Find.Execute("257÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "233÷4=", 2)
$d.Content.Find.Execute("282÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "110÷8=", 2)
$d.Content.Find.Execute("645÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "182÷9=", 2)
$d.Content.Find.Execute("496÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "140÷6=", 2)
$d.Content.Find.Execute("785÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "478÷5=", 2)
$d.Content.Find.Execute("636÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "935÷7=", 2)
$d.Content.Find.Execute("771÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "693÷3=", 2)
$d.Content.Find.Execute("877÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "785÷9=", 2)
$d.Content.Find.Execute("912÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "756÷8=", 2)
$d.Content.Find.Execute("914÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "680÷4=", 2)
